# Applies the commit's changes to the workbook:
#  - Sheet "1.3" ("Controle de Produtos Vendidos") is restructured into two
#    side-by-side mini-tables (a Quantidade/Preço Unitário table and a Preço
#    Unitário/Custo Unitário/Lucro table) followed by an operations header
#    row (Adição / Multiplicação / Subtração / Divisão) with a yellow band.
#    The old B1:E1 title merge is removed.
#  - Sheet "1.4" keeps its original text/values (only the underlying shared
#    string indices shift because of the string-table churn on sheet "1.3";
#    that happens automatically).
#  - The active tab becomes sheet "1.3" (0-based index 2).

$wb = $excel.ActiveWorkbook

$xlLeft   = -4131
$xlRight  = -4152
$xlCenter = -4108

$ws3 = $wb.Worksheets.Item(3)   # "1.3"

# ---------------------------------------------------------------------
# 1) Un-merge the old title band so B1:E1 become independent cells again.
# ---------------------------------------------------------------------
$ws3.Range("B1:E1").UnMerge()

# ---------------------------------------------------------------------
# 2) Write every cell in B1:E13 explicitly (value or blank) so nothing
#    stale survives from the old layout (rows 6-8 previously held the
#    Arroz/Frango rows; those must disappear).
# ---------------------------------------------------------------------

# row -> ordered values for columns B,C,D,E ($null = leave blank)
$rowValues = @{
    1  = @("Controle de Produtos Vendidos", $null, $null, $null)
    2  = @("Produto", "Quantidade", "Preço Unitário", "Total")
    3  = @("Maçãs", 3, 2.5, $null)
    4  = @("Leite", 1, 4, $null)
    5  = @("Pão", 2, 3.25, $null)
    6  = @("Total", $null, $null, $null)
    7  = @($null, $null, $null, $null)
    8  = @("Produto", "Preço Unitário", "Custo Unitário", "Lucro")
    9  = @("Maçãs", 2.5, 1.5, $null)
    10 = @("Leite", 4, 2, $null)
    11 = @("Pão", 3.25, 2.35, $null)
    12 = @($null, $null, $null, $null)
    13 = @("Adição", "Multiplicação", "Subtração", "Divisão")
}

$cols = @("B", "C", "D", "E")

foreach ($r in 1..13) {
    $vals = $rowValues[$r]
    for ($i = 0; $i -lt 4; $i++) {
        $addr = $cols[$i] + $r
        $v = $vals[$i]
        if ($null -eq $v) {
            $ws3.Range($addr).ClearContents()
        } else {
            $ws3.Range($addr).Value = $v
        }
    }
}

# ---------------------------------------------------------------------
# 3) Formatting.
# ---------------------------------------------------------------------

# Whole touched block gets the Calibri 11 "automatic color" font that the
# new table styling uses (a distinct font entry from the default Normal
# style font, even though visually equivalent).
$dataRange = $ws3.Range("B1:E13")
$dataRange.Font.Name = "Calibri"
$dataRange.Font.Size = 11

# Row 1 title: left-aligned label, centered (blank) placeholders D1:E1.
$ws3.Range("B1").HorizontalAlignment = $xlLeft
$ws3.Range("C1").HorizontalAlignment = $xlLeft
$ws3.Range("D1:E1").HorizontalAlignment = $xlCenter

# Header rows (2 and 8) and the text/label column (B) are left-aligned.
$ws3.Range("B2:E2").HorizontalAlignment = $xlLeft
$ws3.Range("B3:B6").HorizontalAlignment = $xlLeft
$ws3.Range("B8:E8").HorizontalAlignment = $xlLeft
$ws3.Range("B9:B11").HorizontalAlignment = $xlLeft

# Numeric columns (Quantidade / Preço Unitário / Custo Unitário) right align.
$ws3.Range("C3:D5").HorizontalAlignment = $xlRight
$ws3.Range("C9:D11").HorizontalAlignment = $xlRight

# Remaining blank-but-styled cells (left aligned, default look).
$leftBlank = @("E3","E4","E5","C6","D6","E6","B7","C7","D7","E7","E9","E10","E11","B12","C12","D12","E12")
foreach ($addr in $leftBlank) {
    $ws3.Range($addr).HorizontalAlignment = $xlLeft
}

# Row 13 operations band: explicit row height + yellow fill + left align.
$ws3.Rows.Item(13).RowHeight = 15
$opRow = $ws3.Range("B13:E13")
$opRow.HorizontalAlignment = $xlLeft
$opRow.Interior.ColorIndex = 6

# ---------------------------------------------------------------------
# 4) Make sheet "1.3" the active tab (workbookView activeTab=2).
# ---------------------------------------------------------------------
$ws3.Activate()
